$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "M1"
$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 0.6435753333333333
$ws.Range("H2").Value2 = 1.930726
$ws.Range("L2").Value2 = 0.5
$ws.Range("M2").Value2 = 2.170377
$ws.Range("N2").Value2 = 4.340754
$ws.Range("O2").Value2 = 0.1015511790371702
$ws.Range("P2").Value2 = 0.07285982038608425
$ws.Range("Q2").Value2 = 1.396801101234
$ws.Range("R2").Value2 = 8.380806607403999
$ws.Range("S2").Value2 = 0.1015511790371702
$ws.Range("T2").Value2 = 0.07285982038608425
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 0.6435753333333333
$ws.Range("H3").Value2 = 1.930726
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 5.061974333333333
$ws.Range("N3").Value2 = 15.185923
$ws.Range("O3").Value2 = 0.2368480046581279
$ws.Range("P3").Value2 = 0.2548966428820674
$ws.Range("Q3").Value2 = 3.257761818899777
$ws.Range("R3").Value2 = 29.319856370098
$ws.Range("S3").Value2 = 0.2368480046581279
$ws.Range("T3").Value2 = 0.2548966428820674
$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 0.6435753333333333
$ws.Range("H4").Value2 = 1.930726
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 3.815520666666667
$ws.Range("N4").Value2 = 11.446562
$ws.Range("O4").Value2 = 0.1785268745202745
$ws.Range("P4").Value2 = 0.1921312406457904
$ws.Range("Q4").Value2 = 2.455574984890222
$ws.Range("R4").Value2 = 22.100174864012
$ws.Range("S4").Value2 = 0.1785268745202745
$ws.Range("T4").Value2 = 0.1921312406457904
$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 0.6435753333333333
$ws.Range("H5").Value2 = 1.930726
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 4.666218666666667
$ws.Range("N5").Value2 = 13.998656
$ws.Range("O5").Value2 = 0.2183307357409577
$ws.Range("P5").Value2 = 0.2349682939430755
$ws.Range("Q5").Value2 = 3.003063233806222
$ws.Range("R5").Value2 = 27.027569104256
$ws.Range("S5").Value2 = 0.2183307357409577
$ws.Range("T5").Value2 = 0.2349682939430755
$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 0.6435753333333333
$ws.Range("H6").Value2 = 1.930726
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 3.288577
$ws.Range("N6").Value2 = 9.865731
$ws.Range("O6").Value2 = 0.1538713650690733
$ws.Range("P6").Value2 = 0.1655968959856798
$ws.Range("Q6").Value2 = 2.116447038967333
$ws.Range("R6").Value2 = 19.048023350706
$ws.Range("S6").Value2 = 0.1538713650690733
$ws.Range("T6").Value2 = 0.1655968959856798
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Wnt7a"
$ws.Range("C7").Value2 = "Fzd5"
$ws.Range("D7").Value2 = "ECs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 0.6435753333333333
$ws.Range("H7").Value2 = 1.930726
$ws.Range("I7").Value2 = 1
$ws.Range("J7").Value2 = 1
$ws.Range("K7").Value2 = 2
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 2.3695805
$ws.Range("N7").Value2 = 4.739161
$ws.Range("O7").Value2 = 0.1108718409743963
$ws.Range("P7").Value2 = 0.07954710615730251
$ws.Range("Q7").Value2 = 1.525003560147667
$ws.Range("R7").Value2 = 9.150021360886
$ws.Range("S7").Value2 = 0.1108718409743963
$ws.Range("T7").Value2 = 0.07954710615730251
